{"js": "// Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line\n// and the \"\u00a9 2020 ... Creative Commons Attribution\" footer line (plus the\n// blank paragraph that used to separate that footer from the page-break\n// paragraph) that followed the bibliography, mirroring a Jekyll site\n// rebuild that dropped the generated page chrome.\n\nconst body = context.document.body;\n\n// Locate the first footer line via search so the edit is anchored to its\n// actual text rather than a hard-coded paragraph index.\nconst anchorText = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst searchResults = body.search(anchorText, { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Anchor paragraph not found: \"' + anchorText + '\"');\n}\n\n// Paragraph containing \"Ver no Jupiter Salvar em pdf Salvar em docx\".\nconst firstParagraph = searchResults.items[0].paragraphs.getFirst();\nfirstParagraph.load(\"text\");\nawait context.sync();\n\n// The paragraph right after it (\"\u00a9 2020 . Contact: ...\").\nconst secondParagraph = firstParagraph.getNext();\nsecondParagraph.load(\"text\");\nawait context.sync();\n\n// The (empty) paragraph right after that one, which also gets removed so\n// the document collapses back to a single blank line before the trailing\n// page-break paragraph.\nconst thirdParagraph = secondParagraph.getNext();\nthirdParagraph.load(\"text\");\nawait context.sync();\n\n// Sanity-check we found the expected text before deleting anything.\nif (secondParagraph.text.indexOf(\"2020\") === -1) {\n  throw new Error('Unexpected paragraph after anchor: \"' + secondParagraph.text + '\"');\n}\nif (thirdParagraph.text !== \"\") {\n  throw new Error('Unexpected trailing paragraph after anchor: \"' + thirdParagraph.text + '\"');\n}\n\n// Delete last-to-first so earlier references stay valid.\nthirdParagraph.delete();\nsecondParagraph.delete();\nfirstParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line\n# and the \"\u00a9 2020 ... Creative Commons Attribution\" footer line (plus the\n# blank paragraph that used to separate that footer from the page-break\n# paragraph) that followed the bibliography, mirroring a Jekyll site\n# rebuild that dropped the generated page chrome.\n\n$d = $word.ActiveDocument\n\n# Anchor on the first footer line's text so the edit targets the right\n# spot even if earlier content in the document shifts around.\n$anchorText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$findRange = $d.Content\n$found = $findRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Anchor paragraph not found: $anchorText\"\n}\n\n# Resolve the paragraph object that contains the found text.\n$startPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Start -le $findRange.Start -and $p.Range.End -ge $findRange.End) {\n        $startPara = $p\n        break\n    }\n}\nif ($null -eq $startPara) {\n    throw \"Could not resolve the anchor paragraph\"\n}\n\n# The footer spans three paragraphs: the anchor line, the \"\u00a9 2020 ...\" line\n# right after it, and the now-orphaned blank paragraph after that.\n$secondPara = $startPara.Next()\n$thirdPara = $secondPara.Next()\n\nif ($secondPara.Range.Text -notlike \"*2020*\") {\n    throw \"Unexpected paragraph after anchor: $($secondPara.Range.Text)\"\n}\nif ($thirdPara.Range.Text.Trim() -ne \"\") {\n    throw \"Unexpected trailing paragraph after anchor: $($thirdPara.Range.Text)\"\n}\n\n$deleteRange = $d.Range($startPara.Range.Start, $thirdPara.Range.End)\n$deleteRange.Delete()\n"}
